# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (Office colours) - used by the notes master
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colours - used by the slide master
#     (and therefore by every slide, since presentation.xml.rels' theme
#     relationship and slideMaster1's theme relationship both point at theme2.xml)
#
# The authored edit swaps the content of the two theme parts, so the slide
# master (and every slide) switches from the "Red Violet" palette over to the
# plain "Office" palette. Re-create that effect through the supported
# PowerPoint object model by rewriting the twelve theme colour slots that
# back the presentation's single (slide-master-facing) theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Office theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as VBA-style BGR long values for the ThemeColor.RGB setter.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeColors[$i - 1]
}
